$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---- Overview sheet ----
$ws1.Range("A1").Value = "File Name"
$ws1.Range("B1").Value = "zh-cn"
$ws1.Range("C1").Value = "de-de"
$ws1.Range("A2").Value = "ffff26b2e6c4-cc09-49e5-abf2-c9e0f58efd83.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("A3").Value = "ffffff93c38c90-4773-4cfd-839e-93171b1bf5b9.md"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"
$ws1.Range("A4").Value = "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("A5").Value = ".localization-config"
$ws1.Range("B5").Value = "Not to be localized"
$ws1.Range("C5").Value = "Not to be localized"
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/97cee173c9365a900001c74c4b53a14c8c2824e5/e2e/5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md", "", "", "ffff26b2e6c4-cc09-49e5-abf2-c9e0f58efd83.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/97cee173c9365a900001c74c4b53a14c8c2824e5/e2e/ffff26b2e6c4-cc09-49e5-abf2-c9e0f58efd83.md", "", "", "ffffff93c38c90-4773-4cfd-839e-93171b1bf5b9.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/97cee173c9365a900001c74c4b53a14c8c2824e5/e2e/ffffff93c38c90-4773-4cfd-839e-93171b1bf5b9.md", "", "", "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/97cee173c9365a900001c74c4b53a14c8c2824e5/.localization-config", "", "", ".localization-config")

# ---- zh-cn sheet ----
$ws2.Range("A1").Value = "Source File Name"
$ws2.Range("B1").Value = "Status"
$ws2.Range("C1").Value = "Latest Handoff File"
$ws2.Range("D1").Value = "Latest Handoff Datetime"
$ws2.Range("E1").Value = "Latest Target File"
$ws2.Range("F1").Value = "Latest Handback File"
$ws2.Range("G1").Value = "Latest Handback DateTime"
$ws2.Range("H1").Value = "Handoff Reason"
$ws2.Range("I1").Value = "Dependency From"
$ws2.Range("A2").Value = "ffff26b2e6c4-cc09-49e5-abf2-c9e0f58efd83.md"
$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-02-22 05:07:48"
$ws2.Range("E2").Value = "48566a70-0a28-4fce-8ad0-9368ac6f1432.md"
$ws2.Range("F2").Value = "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-02-22 05:08:50"
$ws2.Range("H2").Value = "Include"
$ws2.Range("A3").Value = "ffffff93c38c90-4773-4cfd-839e-93171b1bf5b9.md"
$ws2.Range("B3").Value = "Handed back: in sync with en-US"
$ws2.Range("C3").Value = "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-02-22 05:07:48"
$ws2.Range("E3").Value = "48566a70-0a28-4fce-8ad0-9368ac6f1432.md"
$ws2.Range("F3").Value = "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf"
$ws2.Range("G3").Value = "2016-02-22 05:08:50"
$ws2.Range("H3").Value = "Include"
$ws2.Range("A4").Value = "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.6a40e26501acfde0573b2d3789e581706f404b2e.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-02-22 05:13:22"
$ws2.Range("E4").Value = "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md"
$ws2.Range("F4").Value = "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.6a40e26501acfde0573b2d3789e581706f404b2e.zh-cn.xlf"
$ws2.Range("G4").Value = "2016-02-22 05:12:01"
$ws2.Range("H4").Value = "Include"
$ws2.Range("A5").Value = ".localization-config"
$ws2.Range("B5").Value = "Not to be localized"
$ws2.Range("D5").Value = "0001-01-01 00:00:00"
$ws2.Range("G5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Ignored"
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/97cee173c9365a900001c74c4b53a14c8c2824e5/e2e/5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md", "", "", "ffff26b2e6c4-cc09-49e5-abf2-c9e0f58efd83.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/208121a57b04fbad450933099b75d95e0cd24465/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5b344a3f-820b-4a9c-97ba-67ec3b1996e9.6a40e26501acfde0573b2d3789e581706f404b2e.zh-cn.xlf", "", "", "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c92c7e68a367fccb35fae5f8b5afa3bed07683e3/e2e/5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md", "", "", "48566a70-0a28-4fce-8ad0-9368ac6f1432.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c9cd1716b2a23f0d958a304da3bfdf4e0a2c90f8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5b344a3f-820b-4a9c-97ba-67ec3b1996e9.6a40e26501acfde0573b2d3789e581706f404b2e.zh-cn.xlf", "", "", "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/97cee173c9365a900001c74c4b53a14c8c2824e5/e2e/ffff26b2e6c4-cc09-49e5-abf2-c9e0f58efd83.md", "", "", "ffffff93c38c90-4773-4cfd-839e-93171b1bf5b9.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/41e0d1f490bb3de1c1ff71678f7497aae115e44b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf", "", "", "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6bea894f832ed40bef70973eb97310cfde9cc465/e2e/48566a70-0a28-4fce-8ad0-9368ac6f1432.md", "", "", "48566a70-0a28-4fce-8ad0-9368ac6f1432.md")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/af553aee451047ba38d038e57433f37fb989b25f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf", "", "", "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/97cee173c9365a900001c74c4b53a14c8c2824e5/e2e/ffffff93c38c90-4773-4cfd-839e-93171b1bf5b9.md", "", "", "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/41e0d1f490bb3de1c1ff71678f7497aae115e44b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf", "", "", "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.6a40e26501acfde0573b2d3789e581706f404b2e.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6bea894f832ed40bef70973eb97310cfde9cc465/e2e/48566a70-0a28-4fce-8ad0-9368ac6f1432.md", "", "", "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md")
$ws2.Hyperlinks.Add($ws2.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/af553aee451047ba38d038e57433f37fb989b25f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.zh-cn.xlf", "", "", "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.6a40e26501acfde0573b2d3789e581706f404b2e.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/97cee173c9365a900001c74c4b53a14c8c2824e5/.localization-config", "", "", ".localization-config")

# ---- de-de sheet ----
$ws3.Range("A1").Value = "Source File Name"
$ws3.Range("B1").Value = "Status"
$ws3.Range("C1").Value = "Latest Handoff File"
$ws3.Range("D1").Value = "Latest Handoff Datetime"
$ws3.Range("E1").Value = "Latest Target File"
$ws3.Range("F1").Value = "Latest Handback File"
$ws3.Range("G1").Value = "Latest Handback DateTime"
$ws3.Range("H1").Value = "Handoff Reason"
$ws3.Range("I1").Value = "Dependency From"
$ws3.Range("A2").Value = "ffff26b2e6c4-cc09-49e5-abf2-c9e0f58efd83.md"
$ws3.Range("B2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf"
$ws3.Range("D2").Value = "2016-02-22 05:08:03"
$ws3.Range("E2").Value = "48566a70-0a28-4fce-8ad0-9368ac6f1432.md"
$ws3.Range("F2").Value = "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf"
$ws3.Range("G2").Value = "2016-02-22 05:09:15"
$ws3.Range("H2").Value = "Include"
$ws3.Range("A3").Value = "ffffff93c38c90-4773-4cfd-839e-93171b1bf5b9.md"
$ws3.Range("B3").Value = "Handed back: in sync with en-US"
$ws3.Range("C3").Value = "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf"
$ws3.Range("D3").Value = "2016-02-22 05:08:03"
$ws3.Range("E3").Value = "48566a70-0a28-4fce-8ad0-9368ac6f1432.md"
$ws3.Range("F3").Value = "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf"
$ws3.Range("G3").Value = "2016-02-22 05:09:15"
$ws3.Range("H3").Value = "Include"
$ws3.Range("A4").Value = "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.6a40e26501acfde0573b2d3789e581706f404b2e.de-de.xlf"
$ws3.Range("D4").Value = "2016-02-22 05:13:37"
$ws3.Range("E4").Value = "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md"
$ws3.Range("F4").Value = "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.6a40e26501acfde0573b2d3789e581706f404b2e.de-de.xlf"
$ws3.Range("G4").Value = "2016-02-22 05:12:26"
$ws3.Range("H4").Value = "Include"
$ws3.Range("A5").Value = ".localization-config"
$ws3.Range("B5").Value = "Not to be localized"
$ws3.Range("D5").Value = "0001-01-01 00:00:00"
$ws3.Range("G5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Ignored"
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/97cee173c9365a900001c74c4b53a14c8c2824e5/e2e/5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md", "", "", "ffff26b2e6c4-cc09-49e5-abf2-c9e0f58efd83.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f6ec61b14b46a58372c1b324c652cded98a33efd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5b344a3f-820b-4a9c-97ba-67ec3b1996e9.6a40e26501acfde0573b2d3789e581706f404b2e.de-de.xlf", "", "", "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/614ec3bfe2eae03d6e7d4fb90783f36126989b22/e2e/5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md", "", "", "48566a70-0a28-4fce-8ad0-9368ac6f1432.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/03da5e594fd7879685475821bbb4b06fcf6ca929/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5b344a3f-820b-4a9c-97ba-67ec3b1996e9.6a40e26501acfde0573b2d3789e581706f404b2e.de-de.xlf", "", "", "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/97cee173c9365a900001c74c4b53a14c8c2824e5/e2e/ffff26b2e6c4-cc09-49e5-abf2-c9e0f58efd83.md", "", "", "ffffff93c38c90-4773-4cfd-839e-93171b1bf5b9.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2987ccb40f6f9233e2a195117bbe91be24e02410/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf", "", "", "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8199a62641b3d8b27209c7f5e4b08c81229a1ada/e2e/48566a70-0a28-4fce-8ad0-9368ac6f1432.md", "", "", "48566a70-0a28-4fce-8ad0-9368ac6f1432.md")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1e986ddc178d3fee55ebc3d5003046b18fb59e66/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf", "", "", "48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/97cee173c9365a900001c74c4b53a14c8c2824e5/e2e/ffffff93c38c90-4773-4cfd-839e-93171b1bf5b9.md", "", "", "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2987ccb40f6f9233e2a195117bbe91be24e02410/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf", "", "", "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.6a40e26501acfde0573b2d3789e581706f404b2e.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8199a62641b3d8b27209c7f5e4b08c81229a1ada/e2e/48566a70-0a28-4fce-8ad0-9368ac6f1432.md", "", "", "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.md")
$ws3.Hyperlinks.Add($ws3.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1e986ddc178d3fee55ebc3d5003046b18fb59e66/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/48566a70-0a28-4fce-8ad0-9368ac6f1432.171d6ada3299d3719b5b2bdba06903d1e8ad4ae8.de-de.xlf", "", "", "5b344a3f-820b-4a9c-97ba-67ec3b1996e9.6a40e26501acfde0573b2d3789e581706f404b2e.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/97cee173c9365a900001c74c4b53a14c8c2824e5/.localization-config", "", "", ".localization-config")
